# New crime data collected - weekly 121st Precinct CompStat update
# Volume/date header text + full data table refresh (rows 15-30)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Helper source cells for the "0 count" / "***.* pct" placeholder ----
# styled text cells (style 14) already used elsewhere on the sheet, so we
# copy them (value + format together) onto the cells that need to flip
# from a plain number into the dash / N/A placeholder text.
$dashSrc = $ws.Range("C14")   # s=14, t=s -> shared string "0"
$naSrc   = $ws.Range("E14")   # s=14, t=s -> shared string "***.*"

function Set-Dash($addr) {
    $dashSrc.Copy($ws.Range($addr))
}
function Set-NA($addr) {
    $naSrc.Copy($ws.Range($addr))
}

# --- Row 15 : Rape ---------------------------------------------------------
Set-Dash "D15"
Set-NA   "E15"
$ws.Range("L15").Value = -54.545454545454

# --- Row 16 : Robbery -------------------------------------------------------
Set-Dash "C16"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = 11.764705882352

# --- Row 17 : Fel. Assault ---------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -5.263157894736
$ws.Range("I17").Value = 128
$ws.Range("J17").Value = 135
$ws.Range("K17").Value = -5.185185185185
$ws.Range("L17").Value = 60

# --- Row 18 : Burglary -------------------------------------------------------
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = -34.782608695652
$ws.Range("L18").Value = -40

# --- Row 19 : Gr. Larceny ----------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 22.727272727272
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 241
$ws.Range("K19").Value = -13.278008298755
$ws.Range("L19").Value = 28.220858895705

# --- Row 20 : G.L.A. ---------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 59
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = 31.111111111111
$ws.Range("L20").Value = 180.952380952381

# --- Row 21 : TOTAL -----------------------------------------------------------
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -23.529411764705
$ws.Range("F21").Value = 68
$ws.Range("H21").Value = 7.936507936507
$ws.Range("I21").Value = 470
$ws.Range("J21").Value = 507
$ws.Range("K21").Value = -7.297830374753
$ws.Range("L21").Value = 32.022471910112

# --- Row 23 : Housing -----------------------------------------------------------
Set-Dash "G23"
Set-NA   "H23"

# --- Row 24 : Petit Larceny -----------------------------------------------------
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 125
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 27.551020408163
$ws.Range("I24").Value = 746
$ws.Range("J24").Value = 715
$ws.Range("K24").Value = 4.335664335664
$ws.Range("L24").Value = 50.100603621730

# --- Row 25 : Misd. Assault -----------------------------------------------------
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 303
$ws.Range("J25").Value = 274
$ws.Range("K25").Value = 10.583941605839
$ws.Range("L25").Value = 50.746268656716

# --- Row 26 : UCR Rape* ----------------------------------------------------------
$ws.Range("C26").Value = 1
Set-Dash "D26"
Set-NA   "E26"
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 12
$ws.Range("K26").Value = 33.333333333333
$ws.Range("L26").Value = -33.333333333333

# --- Row 27 : Other Sex Crimes -----------------------------------------------------
Set-Dash "C27"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0

# --- Row 28 : Shooting Vic. -----------------------------------------------------
Set-Dash "G28"
Set-NA   "H28"

# --- Row 29 : Shooting Inc. -----------------------------------------------------
Set-Dash "G29"
Set-NA   "H29"

# --- Row 30 : Hate Crimes -----------------------------------------------------
Set-Dash "F30"
$ws.Range("H30").Value = -100
